$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing ParcelId..Ratio
# columns from A:S to B:T.
$ws.Columns("A:A").Insert()

# Excel keeps each destination cell's formatting anchored to its position
# when shifting columns, so the header-row style that the insert left on
# B2:B3 (originally on A2:A3) needs to move back onto the new A2:A3 cells,
# and B2:B3 should fall back to the default (no explicit style).
$ws.Range("B2:B3").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$ws.Range("B2:B3").ClearFormats()

# Populate the new first column with the new id values.
$ws.Range("A2").Value = 44
$ws.Range("A3").Value = 485
